$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-6 from 45170 to 45174
$ws.Range("C2:C6").Value = 45174
